$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 811
$ws.Range("I2").Value = 2137
$ws.Range("J2").Value = 8710
$ws.Range("K2").Value = 37
$ws.Range("L2").Value = 2397
$ws.Range("M2").Value = 137
$ws.Range("N2").Value = 1516
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 29
$ws.Range("Q2").Value = 18
$ws.Range("R2").Value = 114
$ws.Range("S2").Value = 959
$ws.Range("T2").Value = 1555
$ws.Range("U2").Value = 131
$ws.Range("V2").Value = 13529
$ws.Range("W2").Value = 7
$ws.Range("X2").Value = 13443
$ws.Range("Y2").Value = 21
$ws.Range("Z2").Value = 228
$ws.Range("AA2").Value = 102
